$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quarterly")

# Insert two new columns (Month, Quarter) right after column A (Year),
# shifting all the existing quarterly metric columns two places right.
$ws.Columns("B:C").Insert()

# Header row
$ws.Range("B1").Value = "Month"
$ws.Range("C1").Value = "Quarter"

# Column A: replace the combined "Mon YY Qn" label with just the Year,
# stored as text (leading apostrophe) so it matches the shared-string type.
$ws.Range("A2").Value = "'2022"
$ws.Range("A2").Style = "Normal"
$ws.Range("A3").Value = "'2022"
$ws.Range("A3").Style = "Normal"
$ws.Range("A4").Value = "'2022"
$ws.Range("A4").Style = "Normal"
$ws.Range("A5").Value = "'2023"
$ws.Range("A5").Style = "Normal"
$ws.Range("A6").Value = "'2023"
$ws.Range("A6").Style = "Normal"
$ws.Range("A7").Value = "'2023"
$ws.Range("A7").Style = "Normal"
$ws.Range("A8").Value = "'2023"
$ws.Range("A8").Style = "Normal"
$ws.Range("A9").Value = "'2024"
$ws.Range("A9").Style = "Normal"
$ws.Range("A10").Value = "'2024"
$ws.Range("A10").Style = "Normal"
$ws.Range("A11").Value = "'2024"
$ws.Range("A11").Style = "Normal"

# Column B: Month number, kept as zero-padded text
$ws.Range("B2").Value = "'06"
$ws.Range("B2").Style = "Normal"
$ws.Range("B3").Value = "'09"
$ws.Range("B3").Style = "Normal"
$ws.Range("B4").Value = "'12"
$ws.Range("B4").Style = "Normal"
$ws.Range("B5").Value = "'03"
$ws.Range("B5").Style = "Normal"
$ws.Range("B6").Value = "'06"
$ws.Range("B6").Style = "Normal"
$ws.Range("B7").Value = "'09"
$ws.Range("B7").Style = "Normal"
$ws.Range("B8").Value = "'12"
$ws.Range("B8").Style = "Normal"
$ws.Range("B9").Value = "'03"
$ws.Range("B9").Style = "Normal"
$ws.Range("B10").Value = "'06"
$ws.Range("B10").Style = "Normal"
$ws.Range("B11").Value = "'09"
$ws.Range("B11").Style = "Normal"

# Column C: Quarter label
$ws.Range("C2").Value = "Q1"
$ws.Range("C3").Value = "Q2"
$ws.Range("C4").Value = "Q3"
$ws.Range("C5").Value = "Q4"
$ws.Range("C6").Value = "Q1"
$ws.Range("C7").Value = "Q2"
$ws.Range("C8").Value = "Q3"
$ws.Range("C9").Value = "Q4"
$ws.Range("C10").Value = "Q1"
$ws.Range("C11").Value = "Q2"
